$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily practice")

# Row 24: update the Date cell to reflect the vacation + resumed date range
$ws.Cells.Item(24,1).Value = "23rd May, 2025`n(WAS ON VACATION FROM 24th - 27th)`n- 28th May, 2025"

# Row 25 (new): "We need the zero" problem
# (values are written in the same order the author entered them, so new
#  shared-string indices line up with the source workbook)
$ws.Cells.Item(25,1).Value = "28th May, 2025"
$ws.Cells.Item(25,2).Value = "We need the zero"
$ws.Cells.Item(25,3).Value = "https://codeforces.com/problemset/problem/1805/A"
$ws.Cells.Item(25,5).Value = "Established that the XOR of all elements is the answer, but facing trouble when there is a 0 in either the array or the result. Not being able to determine what case should the answer exist and what it should be`nits actually a pretty easy problem, we just need to read the question properly.`nSo, if there are even number of elements, the resultant effectof x is zero since we are xoring x with every element, i.e, there are even number of x's being used and the net effect of that is 0`nSo, if n is even and the xor of the elemetns of the given array is 0, then any x is the answer else there is no answer`nNow, if n is odd, then the xor of all the elements is the result !"
$ws.Cells.Item(25,4).Value = "Easy if you are observant.`nElse, Difficult"
$ws.Rows.Item(25).RowHeight = 115.2

# Row 26 (new): "Prepend and append" problem
$ws.Cells.Item(26,1).Value = "28th May, 2025"
$ws.Cells.Item(26,2).Value = "Prepend and append"
$ws.Cells.Item(26,3).Value = "https://codeforces.com/problemset/problem/1791/C"
$ws.Cells.Item(26,5).Value = "just do the task in a loop and you will eventually reach the answer"
$ws.Cells.Item(26,4).Value = "Easy"
$ws.Rows.Item(26).RowHeight = 28.8

# Activate "Daily practice" tab and scroll/select to mirror the author's final view
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D27").Select()
